$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 42606.882893518516
$ws.Range("A3").NumberFormat = "m/d/yy h:mm"

$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 61
$ws.Range("D3").Value = 37
$ws.Range("E3").Value = 50
$ws.Range("F3").Value = 50
$ws.Range("G3").Value = 4845
$ws.Range("H3").Value = 2595
$ws.Range("I3").Value = 432
$ws.Range("J3").Value = 73
$ws.Range("K3").Value = 45
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = "Noun"
